{"js": "// M\u1eabu 26: remove the leftover \"vnpt.SiteAddress\" placeholder run that\n// followed \"\u0110\u1ecba ch\u1ec9: \" in the \"B\u00ean A\" block, leaving only \"\u0110\u1ecba ch\u1ec9: \".\nconst body = context.document.body;\nconst results = body.search(\"vnpt.SiteAddress\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# M\u1eabu 26: remove the leftover \"vnpt.SiteAddress\" placeholder run that\n# followed \"\u0110\u1ecba ch\u1ec9: \" in the \"B\u00ean A\" block, leaving only \"\u0110\u1ecba ch\u1ec9: \".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"vnpt.SiteAddress\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, \"wdReplaceAll\")\n"}
